# Rename the three logo pictures (two "PearsonLogo" footers + one
# "BTec_Logo-Orange" header) living in this document's first-page /
# default header & footer stories.
#
# wdHeaderFooterIndex:
#   1 = wdHeaderFooterPrimary   (the "default" header/footer)
#   2 = wdHeaderFooterFirstPage (the "first page" header/footer)
#
# InlineShape has no writable .Name, so the standard COM idiom for
# renaming an inline picture is to flip it to a floating Shape, set
# Name there, then flip it back to an InlineShape so the drawing stays
# wp:inline in the saved markup.

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

function Rename-FirstInlinePicture($range, [string]$newName) {
    $inlineShape = $range.InlineShapes.Item(1)
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

# Footer (first page) -> footer1.xml: image1.png -> image2.png
Rename-FirstInlinePicture $section.Footers.Item(2).Range "image2.png"

# Footer (default) -> footer2.xml: image1.png -> image2.png
Rename-FirstInlinePicture $section.Footers.Item(1).Range "image2.png"

# Header (first page) -> header1.xml: image2.jpg -> image1.jpg
Rename-FirstInlinePicture $section.Headers.Item(2).Range "image1.jpg"
